# === Auto-generated edit script ===
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header text (Volume number, report week date range) ----
$ws.Range("A8").Value = "Volume 32   Number  36"
$ws.Range("C9").Value = "Report Covering the Week  9/1/2025  Through  9/7/2025"

# ---- Cells that change category (text <-> number): copy formats from a stable reference cell, then set value ----
$ws.Range("I14").Copy() | Out-Null
$ws.Range("D14").PasteSpecial(-4122) | Out-Null
$ws.Range("D14").Value = 3
$ws.Range("K14").Copy() | Out-Null
$ws.Range("E14").PasteSpecial(-4122) | Out-Null
$ws.Range("E14").Value = -100
$ws.Range("I14").Copy() | Out-Null
$ws.Range("G14").PasteSpecial(-4122) | Out-Null
$ws.Range("G14").Value = 3
$ws.Range("K14").Copy() | Out-Null
$ws.Range("H14").PasteSpecial(-4122) | Out-Null
$ws.Range("H14").Value = -100
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4122) | Out-Null
$ws.Range("D15").Value = "0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("E15").PasteSpecial(-4122) | Out-Null
$ws.Range("E15").Value = "***.*"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4122) | Out-Null
$ws.Range("D20").Value = "0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("E20").PasteSpecial(-4122) | Out-Null
$ws.Range("E20").Value = "***.*"
$ws.Range("I14").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122) | Out-Null
$ws.Range("C22").Value = 1
$ws.Range("I14").Copy() | Out-Null
$ws.Range("D25").PasteSpecial(-4122) | Out-Null
$ws.Range("D25").Value = 1
$ws.Range("K14").Copy() | Out-Null
$ws.Range("E25").PasteSpecial(-4122) | Out-Null
$ws.Range("E25").Value = 0
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4122) | Out-Null
$ws.Range("D27").Value = "0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("E27").PasteSpecial(-4122) | Out-Null
$ws.Range("E27").Value = "***.*"
$ws.Range("I14").Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4122) | Out-Null
$ws.Range("D29").Value = 8
$ws.Range("K14").Copy() | Out-Null
$ws.Range("E29").PasteSpecial(-4122) | Out-Null
$ws.Range("E29").Value = -75
$ws.Range("I14").Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4122) | Out-Null
$ws.Range("D30").Value = 3
$ws.Range("K14").Copy() | Out-Null
$ws.Range("E30").PasteSpecial(-4122) | Out-Null
$ws.Range("E30").Value = -33.333333333333

# ---- Cells with value-only changes (style/category unchanged) ----
$ws.Range("J14").Value = 10
$ws.Range("K14").Value = -70
$ws.Range("L14").Value = -40
$ws.Range("N14").Value = -94.915254237288
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 18
$ws.Range("K15").Value = 12.5
$ws.Range("L15").Value = 50
$ws.Range("M15").Value = 28.571428571428
$ws.Range("N15").Value = -66.666666666666
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 5
$ws.Range("F16").Value = 18
$ws.Range("H16").Value = 28.571428571428
$ws.Range("I16").Value = 124
$ws.Range("J16").Value = 130
$ws.Range("K16").Value = -4.615384615384
$ws.Range("L16").Value = -12.676056338028
$ws.Range("M16").Value = -22.981366459627
$ws.Range("N16").Value = -91.820580474934
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 12
$ws.Range("E17").Value = -33.333333333333
$ws.Range("F17").Value = 35
$ws.Range("G17").Value = 31
$ws.Range("H17").Value = 12.903225806451
$ws.Range("I17").Value = 315
$ws.Range("J17").Value = 264
$ws.Range("K17").Value = 19.318181818181
$ws.Range("L17").Value = 10.13986013986
$ws.Range("M17").Value = 43.835616438356
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = -10
$ws.Range("I18").Value = 102
$ws.Range("J18").Value = 104
$ws.Range("K18").Value = -1.923076923076
$ws.Range("L18").Value = -15.702479338843
$ws.Range("M18").Value = -51.196172248803
$ws.Range("N18").Value = -91.122715404699
$ws.Range("C19").Value = 15
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = 200
$ws.Range("F19").Value = 51
$ws.Range("G19").Value = 23
$ws.Range("H19").Value = 121.739130434783
$ws.Range("I19").Value = 280
$ws.Range("J19").Value = 246
$ws.Range("K19").Value = 13.821138211382
$ws.Range("L19").Value = 1.083032490974
$ws.Range("M19").Value = 40.703517587939
$ws.Range("N19").Value = -29.824561403508
$ws.Range("G20").Value = 12
$ws.Range("H20").Value = -41.666666666666
$ws.Range("I20").Value = 80
$ws.Range("K20").Value = 19.402985074626
$ws.Range("L20").Value = -33.333333333333
$ws.Range("M20").Value = -11.111111111111
$ws.Range("N20").Value = -84.526112185686
$ws.Range("C21").Value = 33
$ws.Range("E21").Value = 26.923076923076
$ws.Range("F21").Value = 122
$ws.Range("G21").Value = 95
$ws.Range("H21").Value = 28.421052631578
$ws.Range("I21").Value = 922
$ws.Range("J21").Value = 837
$ws.Range("K21").Value = 10.155316606929
$ws.Range("L21").Value = -4.257528556593
$ws.Range("M21").Value = 1.991150442477
$ws.Range("N21").Value = -79.851398601398
$ws.Range("E22").Value = 0
$ws.Range("I22").Value = 13
$ws.Range("J22").Value = 13
$ws.Range("L22").Value = -31.578947368421
$ws.Range("M22").Value = -40.90909090909
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = -40
$ws.Range("F23").Value = 10
$ws.Range("G23").Value = 14
$ws.Range("H23").Value = -28.571428571428
$ws.Range("I23").Value = 82
$ws.Range("J23").Value = 83
$ws.Range("K23").Value = -1.204819277108
$ws.Range("L23").Value = -7.865168539325
$ws.Range("M23").Value = 54.716981132075
$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 20
$ws.Range("E24").Value = 10
$ws.Range("F24").Value = 75
$ws.Range("G24").Value = 61
$ws.Range("H24").Value = 22.950819672131
$ws.Range("I24").Value = 597
$ws.Range("J24").Value = 473
$ws.Range("K24").Value = 26.215644820296
$ws.Range("L24").Value = 16.147859922179
$ws.Range("M24").Value = 7.567567567567
$ws.Range("C25").Value = 1
$ws.Range("F25").Value = 6
$ws.Range("G25").Value = 8
$ws.Range("H25").Value = -25
$ws.Range("I25").Value = 61
$ws.Range("J25").Value = 85
$ws.Range("K25").Value = -28.235294117647
$ws.Range("L25").Value = -14.084507042253
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = 14.285714285714
$ws.Range("F26").Value = 23
$ws.Range("G26").Value = 30
$ws.Range("H26").Value = -23.333333333333
$ws.Range("I26").Value = 306
$ws.Range("J26").Value = 349
$ws.Range("K26").Value = -12.320916905444
$ws.Range("L26").Value = -22.53164556962
$ws.Range("M26").Value = -45.648312611012
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 23
$ws.Range("K27").Value = -23.333333333333
$ws.Range("L27").Value = 15
$ws.Range("C28").Value = 5
$ws.Range("E28").Value = 150
$ws.Range("F28").Value = 14
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = 133.333333333333
$ws.Range("I28").Value = 58
$ws.Range("J28").Value = 40
$ws.Range("K28").Value = 45
$ws.Range("L28").Value = 11.538461538461
$ws.Range("C29").Value = 2
$ws.Range("F29").Value = 3
$ws.Range("G29").Value = 10
$ws.Range("H29").Value = -70
$ws.Range("I29").Value = 14
$ws.Range("J29").Value = 30
$ws.Range("K29").Value = -53.333333333333
$ws.Range("L29").Value = -22.222222222222
$ws.Range("M29").Value = -73.076923076923
$ws.Range("N29").Value = -92.783505154639
$ws.Range("C30").Value = 2
$ws.Range("F30").Value = 3
$ws.Range("G30").Value = 5
$ws.Range("H30").Value = -40
$ws.Range("I30").Value = 11
$ws.Range("J30").Value = 21
$ws.Range("K30").Value = -47.619047619047
$ws.Range("L30").Value = -31.25
$ws.Range("M30").Value = -71.794871794871
$ws.Range("N30").Value = -93.714285714285
$ws.Range("G31").Value = 2
$ws.Range("J31").Value = 3
$ws.Range("K31").Value = -33.333333333333

$excel.CutCopyMode = 0
